# ---------------------------------------------------------------------------
# READY_FOR_SPELLING becomes USE_SPELLING: extend the recap sheet with four
# more dialogue blocks (DIALOGUE 9..12), mechanically cycling the generation
# on S (S = 0, 2, 4, 8, 16, 32) the same way DIALOGUE 1..8 already do.
#
# Each dialogue block is 8 rows tall:
#   header row  -> dialogue label repeated in A/F/K/P/U/Z, "S = 0" in C,
#                  and "S = 2"/"S = 4"/"S = 8"/"S = 16"/"S = 32" in H/M/R/W/AB
#   5 data rows -> numeric pairs (mean, n) under each header pair of columns
#   2 blank rows -> separator before the next block
#
# The new blocks simply continue the existing layout starting right after
# the DIALOGUE 8 block (which ends at row 66), at rows 69, 77, 85 and 93.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new DIALOGUE 9-12 blocks (rows 69-98) mirroring the existing
# DIALOGUE 1-8 layout: a header row (dialogue label / "S = 0" / "S = 2" / "S = 4" / "S = 8" / "S = 16" / "S = 32")
# followed by five data rows, repeated for each of the four new dialogues, separated by two blank rows.

# --- DIALOGUE 9 header (row 69) ---
$ws.Range("A69").Value = "DIALOGUE 9"
$ws.Range("C69").Value = "S = 0"
$ws.Range("F69").Value = "DIALOGUE 9"
$ws.Range("H69").Value = "S = 2"
$ws.Range("K69").Value = "DIALOGUE 9"
$ws.Range("M69").Value = "S = 4"
$ws.Range("P69").Value = "DIALOGUE 9"
$ws.Range("R69").Value = "S = 8"
$ws.Range("U69").Value = "DIALOGUE 9"
$ws.Range("W69").Value = "S = 16"
$ws.Range("Z69").Value = "DIALOGUE 9"
$ws.Range("AB69").Value = "S = 32"

$ws.Range("A70").Value = 1.5
$ws.Range("B70").Value = 27
$ws.Range("C70").Value = 0.88
$ws.Range("D70").Value = 23
$ws.Range("F70").Value = 1.88
$ws.Range("G70").Value = 30
$ws.Range("H70").Value = 0.71
$ws.Range("I70").Value = 23
$ws.Range("K70").Value = 1.71
$ws.Range("L70").Value = 36
$ws.Range("M70").Value = 0.79
$ws.Range("N70").Value = 27
$ws.Range("P70").Value = 2.59
$ws.Range("Q70").Value = 44
$ws.Range("R70").Value = 0.88
$ws.Range("S70").Value = 22
$ws.Range("U70").Value = 1.93
$ws.Range("V70").Value = 27
$ws.Range("W70").Value = 1
$ws.Range("X70").Value = 18
$ws.Range("Z70").Value = 1.38
$ws.Range("AA70").Value = 47
$ws.Range("AB70").Value = 0.31
$ws.Range("AC70").Value = 55

$ws.Range("A71").Value = 1.47
$ws.Range("B71").Value = 25
$ws.Range("C71").Value = 0.85
$ws.Range("D71").Value = 22
$ws.Range("F71").Value = 1.7
$ws.Range("G71").Value = 34
$ws.Range("H71").Value = 0.52
$ws.Range("I71").Value = 29
$ws.Range("K71").Value = 1.71
$ws.Range("L71").Value = 29
$ws.Range("M71").Value = 0.85
$ws.Range("N71").Value = 22
$ws.Range("P71").Value = 1.88
$ws.Range("Q71").Value = 32
$ws.Range("R71").Value = 0.85
$ws.Range("S71").Value = 22
$ws.Range("U71").Value = 1.47
$ws.Range("V71").Value = 28
$ws.Range("W71").Value = 0.73
$ws.Range("X71").Value = 25
$ws.Range("Z71").Value = 1.76
$ws.Range("AA71").Value = 30
$ws.Range("AB71").Value = 0.85
$ws.Range("AC71").Value = 22

$ws.Range("A72").Value = 2.33
$ws.Range("B72").Value = 14
$ws.Range("C72").Value = 0.6899999999999999
$ws.Range("D72").Value = 13
$ws.Range("F72").Value = 5
$ws.Range("G72").Value = 35
$ws.Range("H72").Value = 0.6
$ws.Range("I72").Value = 15
$ws.Range("K72").Value = 4.62
$ws.Range("L72").Value = 37
$ws.Range("M72").Value = 0.53
$ws.Range("N72").Value = 17
$ws.Range("P72").Value = 2.33
$ws.Range("Q72").Value = 14
$ws.Range("R72").Value = 0.82
$ws.Range("S72").Value = 12
$ws.Range("U72").Value = 5
$ws.Range("V72").Value = 30
$ws.Range("W72").Value = 0.6899999999999999
$ws.Range("X72").Value = 13
$ws.Range("Z72").Value = 2.4
$ws.Range("AA72").Value = 12
$ws.Range("AB72").Value = 0.82
$ws.Range("AC72").Value = 11

$ws.Range("A73").Value = 1.35
$ws.Range("B73").Value = 69
$ws.Range("C73").Value = 0.91
$ws.Range("D73").Value = 65
$ws.Range("F73").Value = 1.71
$ws.Range("G73").Value = 96
$ws.Range("H73").Value = 0.83
$ws.Range("I73").Value = 72
$ws.Range("K73").Value = 1.69
$ws.Range("L73").Value = 98
$ws.Range("M73").Value = 0.91
$ws.Range("N73").Value = 72
$ws.Range("P73").Value = 1.77
$ws.Range("Q73").Value = 92
$ws.Range("R73").Value = 1
$ws.Range("S73").Value = 64
$ws.Range("U73").Value = 1.62
$ws.Range("V73").Value = 91
$ws.Range("W73").Value = 0.91
$ws.Range("X73").Value = 70
$ws.Range("Z73").Value = 1.81
$ws.Range("AA73").Value = 96
$ws.Range("AB73").Value = 0.87
$ws.Range("AC73").Value = 68

$ws.Range("A74").Value = 2.09
$ws.Range("B74").Value = 23
$ws.Range("C74").Value = 0.87
$ws.Range("D74").Value = 17
$ws.Range("F74").Value = 1.8
$ws.Range("G74").Value = 27
$ws.Range("H74").Value = 0.68
$ws.Range("I74").Value = 23
$ws.Range("K74").Value = 1.82
$ws.Range("L74").Value = 20
$ws.Range("M74").Value = 0.87
$ws.Range("N74").Value = 17
$ws.Range("P74").Value = 1.92
$ws.Range("Q74").Value = 23
$ws.Range("R74").Value = 0.87
$ws.Range("S74").Value = 18
$ws.Range("U74").Value = 2.2
$ws.Range("V74").Value = 22
$ws.Range("W74").Value = 0.87
$ws.Range("X74").Value = 16
$ws.Range("Z74").Value = 1.93
$ws.Range("AA74").Value = 27
$ws.Range("AB74").Value = 0.87
$ws.Range("AC74").Value = 20

# --- DIALOGUE 10 header (row 77) ---
$ws.Range("A77").Value = "DIALOGUE 10"
$ws.Range("C77").Value = "S = 0"
$ws.Range("F77").Value = "DIALOGUE 10"
$ws.Range("H77").Value = "S = 2"
$ws.Range("K77").Value = "DIALOGUE 10"
$ws.Range("M77").Value = "S = 4"
$ws.Range("P77").Value = "DIALOGUE 10"
$ws.Range("R77").Value = "S = 8"
$ws.Range("U77").Value = "DIALOGUE 10"
$ws.Range("W77").Value = "S = 16"
$ws.Range("Z77").Value = "DIALOGUE 10"
$ws.Range("AB77").Value = "S = 32"

$ws.Range("A78").Value = 1.81
$ws.Range("B78").Value = 29
$ws.Range("C78").Value = 0.79
$ws.Range("D78").Value = 22
$ws.Range("F78").Value = 1.81
$ws.Range("G78").Value = 38
$ws.Range("H78").Value = 0.71
$ws.Range("I78").Value = 28
$ws.Range("K78").Value = 1.53
$ws.Range("L78").Value = 23
$ws.Range("M78").Value = 0.88
$ws.Range("N78").Value = 20
$ws.Range("P78").Value = 1.8
$ws.Range("Q78").Value = 27
$ws.Range("R78").Value = 0.88
$ws.Range("S78").Value = 20
$ws.Range("U78").Value = 1.44
$ws.Range("V78").Value = 23
$ws.Range("W78").Value = 0.88
$ws.Range("X78").Value = 21
$ws.Range("Z78").Value = 1.9
$ws.Range("AA78").Value = 38
$ws.Range("AB78").Value = 0.79
$ws.Range("AC78").Value = 26

$ws.Range("A79").Value = 1.81
$ws.Range("B79").Value = 29
$ws.Range("C79").Value = 0.85
$ws.Range("D79").Value = 21
$ws.Range("F79").Value = 1.59
$ws.Range("G79").Value = 43
$ws.Range("H79").Value = 0.41
$ws.Range("I79").Value = 39
$ws.Range("K79").Value = 1.79
$ws.Range("L79").Value = 34
$ws.Range("M79").Value = 0.65
$ws.Range("N79").Value = 26
$ws.Range("P79").Value = 1.71
$ws.Range("Q79").Value = 29
$ws.Range("R79").Value = 0.85
$ws.Range("S79").Value = 22
$ws.Range("U79").Value = 1.71
$ws.Range("V79").Value = 29
$ws.Range("W79").Value = 0.85
$ws.Range("X79").Value = 22
$ws.Range("Z79").Value = 2
$ws.Range("AA79").Value = 34
$ws.Range("AB79").Value = 0.85
$ws.Range("AC79").Value = 22

$ws.Range("A80").Value = 8.199999999999999
$ws.Range("B80").Value = 41
$ws.Range("C80").Value = 0.82
$ws.Range("D80").Value = 11
$ws.Range("F80").Value = 4.33
$ws.Range("G80").Value = 39
$ws.Range("H80").Value = 0.47
$ws.Range("I80").Value = 19
$ws.Range("K80").Value = 5
$ws.Range("L80").Value = 35
$ws.Range("M80").Value = 0.6
$ws.Range("N80").Value = 15
$ws.Range("P80").Value = 3.67
$ws.Range("Q80").Value = 33
$ws.Range("R80").Value = 0.47
$ws.Range("S80").Value = 19
$ws.Range("U80").Value = 2.2
$ws.Range("V80").Value = 11
$ws.Range("W80").Value = 0.82
$ws.Range("X80").Value = 11
$ws.Range("Z80").Value = 5.8
$ws.Range("AA80").Value = 29
$ws.Range("AB80").Value = 0.82
$ws.Range("AC80").Value = 11

$ws.Range("A81").Value = 1.38
$ws.Range("B81").Value = 69
$ws.Range("C81").Value = 1
$ws.Range("D81").Value = 62
$ws.Range("F81").Value = 1.76
$ws.Range("G81").Value = 97
$ws.Range("H81").Value = 0.87
$ws.Range("I81").Value = 70
$ws.Range("K81").Value = 1.68
$ws.Range("L81").Value = 84
$ws.Range("M81").Value = 0.95
$ws.Range("N81").Value = 63
$ws.Range("P81").Value = 1.39
$ws.Range("Q81").Value = 82
$ws.Range("R81").Value = 0.95
$ws.Range("S81").Value = 72
$ws.Range("U81").Value = 1.62
$ws.Range("V81").Value = 104
$ws.Range("W81").Value = 0.95
$ws.Range("X81").Value = 77
$ws.Range("Z81").Value = 1.69
$ws.Range("AA81").Value = 105
$ws.Range("AB81").Value = 0.91
$ws.Range("AC81").Value = 76

$ws.Range("A82").Value = 2.08
$ws.Range("B82").Value = 27
$ws.Range("C82").Value = 0.76
$ws.Range("D82").Value = 20
$ws.Range("F82").Value = 1.71
$ws.Range("G82").Value = 29
$ws.Range("H82").Value = 0.76
$ws.Range("I82").Value = 24
$ws.Range("K82").Value = 2.07
$ws.Range("L82").Value = 29
$ws.Range("M82").Value = 0.87
$ws.Range("N82").Value = 20
$ws.Range("P82").Value = 1.93
$ws.Range("Q82").Value = 29
$ws.Range("R82").Value = 0.87
$ws.Range("S82").Value = 21
$ws.Range("U82").Value = 1.57
$ws.Range("V82").Value = 22
$ws.Range("W82").Value = 0.87
$ws.Range("X82").Value = 20
$ws.Range("Z82").Value = 1.92
$ws.Range("AA82").Value = 25
$ws.Range("AB82").Value = 0.76
$ws.Range("AC82").Value = 20

# --- DIALOGUE 11 header (row 85) ---
$ws.Range("A85").Value = "DIALOGUE 11"
$ws.Range("C85").Value = "S = 0"
$ws.Range("F85").Value = "DIALOGUE 11"
$ws.Range("H85").Value = "S = 2"
$ws.Range("K85").Value = "DIALOGUE 11"
$ws.Range("M85").Value = "S = 4"
$ws.Range("P85").Value = "DIALOGUE 11"
$ws.Range("R85").Value = "S = 8"
$ws.Range("U85").Value = "DIALOGUE 11"
$ws.Range("W85").Value = "S = 16"
$ws.Range("Z85").Value = "DIALOGUE 11"
$ws.Range("AB85").Value = "S = 32"

$ws.Range("A86").Value = 1.94
$ws.Range("B86").Value = 33
$ws.Range("C86").Value = 0.88
$ws.Range("D86").Value = 22
$ws.Range("F86").Value = 1.4
$ws.Range("G86").Value = 28
$ws.Range("H86").Value = 0.71
$ws.Range("I86").Value = 27
$ws.Range("K86").Value = 1.63
$ws.Range("L86").Value = 31
$ws.Range("M86").Value = 0.88
$ws.Range("N86").Value = 24
$ws.Range("P86").Value = 1.89
$ws.Range("Q86").Value = 34
$ws.Range("R86").Value = 0.88
$ws.Range("S86").Value = 23
$ws.Range("U86").Value = 2.18
$ws.Range("V86").Value = 37
$ws.Range("W86").Value = 0.88
$ws.Range("X86").Value = 22
$ws.Range("Z86").Value = 1.85
$ws.Range("AA86").Value = 24
$ws.Range("AB86").Value = 0.88
$ws.Range("AC86").Value = 18

$ws.Range("A87").Value = 1.72
$ws.Range("B87").Value = 31
$ws.Range("C87").Value = 0.85
$ws.Range("D87").Value = 23
$ws.Range("F87").Value = 1.7
$ws.Range("G87").Value = 34
$ws.Range("H87").Value = 0.65
$ws.Range("I87").Value = 27
$ws.Range("K87").Value = 1.65
$ws.Range("L87").Value = 33
$ws.Range("M87").Value = 0.85
$ws.Range("N87").Value = 25
$ws.Range("P87").Value = 1.81
$ws.Range("Q87").Value = 29
$ws.Range("R87").Value = 0.85
$ws.Range("S87").Value = 21
$ws.Range("U87").Value = 1.71
$ws.Range("V87").Value = 29
$ws.Range("W87").Value = 1
$ws.Range("X87").Value = 21
$ws.Range("Z87").Value = 1.88
$ws.Range("AA87").Value = 32
$ws.Range("AB87").Value = 0.85
$ws.Range("AC87").Value = 22

$ws.Range("A88").Value = 5.4
$ws.Range("B88").Value = 27
$ws.Range("C88").Value = 0.82
$ws.Range("D88").Value = 11
$ws.Range("F88").Value = 1.73
$ws.Range("G88").Value = 19
$ws.Range("H88").Value = 0.39
$ws.Range("I88").Value = 23
$ws.Range("K88").Value = 5
$ws.Range("L88").Value = 35
$ws.Range("M88").Value = 0.6
$ws.Range("N88").Value = 15
$ws.Range("P88").Value = 2.5
$ws.Range("Q88").Value = 15
$ws.Range("R88").Value = 0.6899999999999999
$ws.Range("S88").Value = 13
$ws.Range("U88").Value = 2.6
$ws.Range("V88").Value = 13
$ws.Range("W88").Value = 0.82
$ws.Range("X88").Value = 11
$ws.Range("Z88").Value = 5.4
$ws.Range("AA88").Value = 27
$ws.Range("AB88").Value = 0.82
$ws.Range("AC88").Value = 11

$ws.Range("A89").Value = 1.73
$ws.Range("B89").Value = 95
$ws.Range("C89").Value = 0.76
$ws.Range("D89").Value = 73
$ws.Range("F89").Value = 1.35
$ws.Range("G89").Value = 69
$ws.Range("H89").Value = 0.8
$ws.Range("I89").Value = 68
$ws.Range("K89").Value = 1.78
$ws.Range("L89").Value = 87
$ws.Range("M89").Value = 0.95
$ws.Range("N89").Value = 62
$ws.Range("P89").Value = 1.36
$ws.Range("Q89").Value = 75
$ws.Range("R89").Value = 0.91
$ws.Range("S89").Value = 69
$ws.Range("U89").Value = 1.63
$ws.Range("V89").Value = 88
$ws.Range("W89").Value = 0.8
$ws.Range("X89").Value = 71
$ws.Range("Z89").Value = 1.39
$ws.Range("AA89").Value = 82
$ws.Range("AB89").Value = 0.87
$ws.Range("AC89").Value = 74

$ws.Range("A90").Value = 2
$ws.Range("B90").Value = 28
$ws.Range("C90").Value = 0.76
$ws.Range("D90").Value = 21
$ws.Range("F90").Value = 1.64
$ws.Range("G90").Value = 23
$ws.Range("H90").Value = 0.57
$ws.Range("I90").Value = 24
$ws.Range("K90").Value = 1.85
$ws.Range("L90").Value = 24
$ws.Range("M90").Value = 0.68
$ws.Range("N90").Value = 21
$ws.Range("P90").Value = 2.08
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = 0.87
$ws.Range("S90").Value = 18
$ws.Range("U90").Value = 2
$ws.Range("V90").Value = 28
$ws.Range("W90").Value = 0.57
$ws.Range("X90").Value = 24
$ws.Range("Z90").Value = 2.18
$ws.Range("AA90").Value = 24
$ws.Range("AB90").Value = 0.76
$ws.Range("AC90").Value = 18

# --- DIALOGUE 12 header (row 93) ---
$ws.Range("A93").Value = "DIALOGUE 12"
$ws.Range("C93").Value = "S = 0"
$ws.Range("F93").Value = "DIALOGUE 12"
$ws.Range("H93").Value = "S = 2"
$ws.Range("K93").Value = "DIALOGUE 12"
$ws.Range("M93").Value = "S = 4"
$ws.Range("P93").Value = "DIALOGUE 12"
$ws.Range("R93").Value = "S = 8"
$ws.Range("U93").Value = "DIALOGUE 12"
$ws.Range("W93").Value = "S = 16"
$ws.Range("Z93").Value = "DIALOGUE 12"
$ws.Range("AB93").Value = "S = 32"

$ws.Range("A94").Value = 1.81
$ws.Range("B94").Value = 29
$ws.Range("C94").Value = 0.88
$ws.Range("D94").Value = 21
$ws.Range("F94").Value = 1.39
$ws.Range("G94").Value = 25
$ws.Range("H94").Value = 0.6
$ws.Range("I94").Value = 27
$ws.Range("K94").Value = 1.42
$ws.Range("L94").Value = 27
$ws.Range("M94").Value = 0.88
$ws.Range("N94").Value = 24
$ws.Range("P94").Value = 1.76
$ws.Range("Q94").Value = 30
$ws.Range("R94").Value = 0.88
$ws.Range("S94").Value = 22
$ws.Range("U94").Value = 1.64
$ws.Range("V94").Value = 36
$ws.Range("W94").Value = 0.71
$ws.Range("X94").Value = 29
$ws.Range("Z94").Value = 1.76
$ws.Range("AA94").Value = 30
$ws.Range("AB94").Value = 0.88
$ws.Range("AC94").Value = 22

$ws.Range("A95").Value = 1.68
$ws.Range("B95").Value = 32
$ws.Range("C95").Value = 0.73
$ws.Range("D95").Value = 25
$ws.Range("F95").Value = 1.68
$ws.Range("G95").Value = 37
$ws.Range("H95").Value = 0.48
$ws.Range("I95").Value = 32
$ws.Range("K95").Value = 1.47
$ws.Range("L95").Value = 28
$ws.Range("M95").Value = 0.65
$ws.Range("N95").Value = 26
$ws.Range("P95").Value = 1.57
$ws.Range("Q95").Value = 33
$ws.Range("R95").Value = 0.58
$ws.Range("S95").Value = 29
$ws.Range("U95").Value = 1.47
$ws.Range("V95").Value = 25
$ws.Range("W95").Value = 0.85
$ws.Range("X95").Value = 22
$ws.Range("Z95").Value = 1.65
$ws.Range("AA95").Value = 33
$ws.Range("AB95").Value = 0.65
$ws.Range("AC95").Value = 27

$ws.Range("A96").Value = 2
$ws.Range("B96").Value = 12
$ws.Range("C96").Value = 0.6899999999999999
$ws.Range("D96").Value = 13
$ws.Range("F96").Value = 4.5
$ws.Range("G96").Value = 27
$ws.Range("H96").Value = 0.6899999999999999
$ws.Range("I96").Value = 13
$ws.Range("K96").Value = 2.6
$ws.Range("L96").Value = 13
$ws.Range("M96").Value = 0.82
$ws.Range("N96").Value = 11
$ws.Range("P96").Value = 2.6
$ws.Range("Q96").Value = 13
$ws.Range("R96").Value = 0.82
$ws.Range("S96").Value = 11
$ws.Range("U96").Value = 5.6
$ws.Range("V96").Value = 28
$ws.Range("W96").Value = 0.82
$ws.Range("X96").Value = 11
$ws.Range("Z96").Value = 6.6
$ws.Range("AA96").Value = 33
$ws.Range("AB96").Value = 0.82
$ws.Range("AC96").Value = 11

$ws.Range("A97").Value = 1.37
$ws.Range("B97").Value = 78
$ws.Range("C97").Value = 0.95
$ws.Range("D97").Value = 70
$ws.Range("F97").Value = 1.58
$ws.Range("G97").Value = 104
$ws.Range("H97").Value = 0.95
$ws.Range("I97").Value = 79
$ws.Range("K97").Value = 1.76
$ws.Range("L97").Value = 97
$ws.Range("M97").Value = 0.95
$ws.Range("N97").Value = 68
$ws.Range("P97").Value = 1.78
$ws.Range("Q97").Value = 89
$ws.Range("R97").Value = 0.95
$ws.Range("S97").Value = 63
$ws.Range("U97").Value = 1.58
$ws.Range("V97").Value = 144
$ws.Range("W97").Value = 0.83
$ws.Range("X97").Value = 107
$ws.Range("Z97").Value = 1.94
$ws.Range("AA97").Value = 95
$ws.Range("AB97").Value = 1
$ws.Range("AC97").Value = 61

$ws.Range("A98").Value = 1.59
$ws.Range("B98").Value = 27
$ws.Range("C98").Value = 0.62
$ws.Range("D98").Value = 26
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 22
$ws.Range("H98").Value = 0.62
$ws.Range("I98").Value = 20
$ws.Range("K98").Value = 2.5
$ws.Range("L98").Value = 20
$ws.Range("M98").Value = 0.87
$ws.Range("N98").Value = 14
$ws.Range("P98").Value = 1.86
$ws.Range("Q98").Value = 26
$ws.Range("R98").Value = 0.87
$ws.Range("S98").Value = 20
$ws.Range("U98").Value = 1.92
$ws.Range("V98").Value = 23
$ws.Range("W98").Value = 0.87
$ws.Range("X98").Value = 18
$ws.Range("Z98").Value = 1.62
$ws.Range("AA98").Value = 21
$ws.Range("AB98").Value = 0.76
$ws.Range("AC98").Value = 20
